# Add a new "Deduplicated" worksheet before the existing data sheet,
# rename the existing sheet to "All codes", and populate the new sheet
# with one row per SNOMED code / missing-language combination.

$wb = $excel.ActiveWorkbook

$allCodes = $wb.Worksheets.Item(1)
$allCodes.Name = "All codes"

$dedup = $wb.Worksheets.Add($allCodes)
$dedup.Name = "Deduplicated"

# Header row
$dedup.Cells.Item(1,1).Value = "SNOMED Code"
$dedup.Cells.Item(1,2).Value = "Display"
$dedup.Cells.Item(1,3).Value = "Missing language"

$dedupData = @(
     @("10179008", "Invasive", "NL")
    ,@("10179008", "Invasive", "FR")
    ,@("14803004", "Transitory", "NL")
    ,@("14803004", "Transitory", "FR")
    ,@("22762002", "Non-invasive", "NL")
    ,@("22762002", "Non-invasive", "FR")
    ,@("26643006", "Oral route", "NL")
    ,@("26643006", "Oral route", "FR")
    ,@("33553000", "Widowed", "NL")
    ,@("33553000", "Widowed", "FR")
    ,@("46713006", "Nasal route", "NL")
    ,@("46713006", "Nasal route", "FR")
    ,@("47625008", "Intravenous route", "NL")
    ,@("47625008", "Intravenous route", "FR")
    ,@("69620002", "Daily", "NL")
    ,@("69620002", "Daily", "FR")
    ,@("74964007", "Other", "NL")
    ,@("74964007", "Other", "FR")
    ,@("84638005", "Occasional", "NL")
    ,@("84638005", "Occasional", "FR")
    ,@("89925002", "Canceled", "NL")
    ,@("89925002", "Canceled", "FR")
    ,@("125681006", "Single person", "NL")
    ,@("125681006", "Single person", "FR")
    ,@("160504008", "Marital state unknown", "FR")
    ,@("160708008", "Stairs in house", "FR")
    ,@("183446009", "Admission to cardiac intensive care unit", "NL")
    ,@("183446009", "Admission to cardiac intensive care unit", "FR")
    ,@("183966005", "Drug treatment not indicated", "FR")
    ,@("184003006", "Patient requested treatment", "FR")
    ,@("186065003", "Power of attorney medical report", "NL")
    ,@("186065003", "Power of attorney medical report", "FR")
    ,@("224152001", "Keeps pets", "FR")
    ,@("224670004", "Houseboat", "NL")
    ,@("224670004", "Houseboat", "FR")
    ,@("248160001", "Undressed", "FR")
    ,@("255545003", "Definite", "NL")
    ,@("255545003", "Definite", "FR")
    ,@("257564005", "Apartment", "NL")
    ,@("257564005", "Apartment", "FR")
    ,@("257630004", "House", "NL")
    ,@("257630004", "House", "FR")
    ,@("257641002", "Mobile home", "NL")
    ,@("257641002", "Mobile home", "FR")
    ,@("260413007", "None", "NL")
    ,@("260413007", "None", "FR")
    ,@("261665006", "Unknown", "NL")
    ,@("261665006", "Unknown", "FR")
    ,@("263718001", "Complication", "NL")
    ,@("263718001", "Complication", "FR")
    ,@("274512008", "Medication discontinued", "FR")
    ,@("276727009", "Null", "NL")
    ,@("276727009", "Null", "FR")
    ,@("277132007", "Therapeutic procedure", "FR")
    ,@("371911009", "Measurement of blood pressure using cuff method", "FR")
    ,@("385643006", "To be done", "NL")
    ,@("385643006", "To be done", "FR")
    ,@("385655000", "Suspended", "NL")
    ,@("385655000", "Suspended", "FR")
    ,@("405613005", "Planned procedure", "FR")
    ,@("406149000", "Medication declined", "FR")
    ,@("408366001", "Unable to use medication", "FR")
    ,@("410521004", "Not to be done", "NL")
    ,@("410521004", "Not to be done", "FR")
    ,@("417284009", "Current drug user", "NL")
    ,@("417284009", "Current drug user", "FR")
    ,@("423155007", "Housing contains structural barriers to movement", "FR")
    ,@("423527000", "Unsafe floor covering", "FR")
    ,@("424415008", "Dirty living conditions", "FR")
    ,@("424661000", "Cluttered living space", "FR")
    ,@("424948003", "Obstructed means of residential entrance", "FR")
    ,@("447694001", "Respiratory tract route", "NL")
    ,@("447694001", "Respiratory tract route", "FR")
    ,@("465153004", "Stairlift, chair", "NL")
    ,@("465153004", "Stairlift, chair", "FR")
    ,@("704417003", "At risk of medication side effect", "FR")
    ,@("705046008", "Deficit of medication supply", "FR")
    ,@("705390009", "Home and premises assistive furnishing/adaptation", "NL")
    ,@("705390009", "Home and premises assistive furnishing/adaptation", "FR")
    ,@("705401004", "Assistive bed", "NL")
    ,@("705401004", "Assistive bed", "FR")
    ,@("713914004", "User of smokeless tobacco", "FR")
    ,@("715758005", "Access to residence by stairs", "FR")
    ,@("715881003", "Carbapenemase producing Enterobacteriaceae carrier", "FR")
    ,@("722499006", "Electronic cigarette user", "FR")
    ,@("768832004", "Case manager", "NL")
    ,@("768832004", "Case manager", "FR")
    ,@("783261004", "Lifetime non-drinker of alcohol", "FR")
    ,@("785889008", "Nicotine-filled electronic cigarette user", "FR")
    ,@("786063001", "Non-nicotine-filled electronic cigarette user", "FR")
    ,@(32570681000036100, "Indeterminate sex", "FR")
)

$r = 2
foreach ($row in $dedupData) {
    $dedup.Cells.Item($r,1).Value = $row[0]
    $dedup.Cells.Item($r,2).Value = $row[1]
    $dedup.Cells.Item($r,3).Value = $row[2]
    $r = $r + 1
}

# Column widths on the new "Deduplicated" sheet
$dedup.Columns.Item(1).ColumnWidth = 9.285714285714286
$dedup.Columns.Item(2).ColumnWidth = 38.57142857142857
$dedup.Columns.Item(4).ColumnWidth = 10.714285714285714

# Column width tweak on "All codes" (A: 12 -> 10)
$allCodes.Columns.Item(1).ColumnWidth = 9.285714285714286

# Selections / active sheet
$dedup.Range("B25").Select()
$dedup.Activate()

$allCodes.Range("A77:B179").Select()
$allCodes.Range("A77").Activate()
